$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new sheet "2022-Q1" right after "2021-Q4" (i.e. right
#    before "总计"), holding the fund-holding breakdown for that quarter.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $src)
$ws.Name = "2022-Q1"

# Header row - copy the sibling sheet's header formatting (bold + border
# style), then overwrite the text.
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Column A is the 0-based row index, styled the same way as the header.
$src.Range("A2:A3").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# Data cells (B:G) are plain, unstyled text cells - even the numeric-
# looking ones ("15.30", "0.6916", ...) are literal text, not numbers.
# Force text number-format while assigning so Excel doesn't coerce them,
# then clear the format again so no style index lingers on the cells.
$dataRange = $ws.Range("B2:G4")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "000727"
$ws.Range("C2").Value = "融通健康产业灵活配置混合A"
$ws.Range("D2").Value = "15.30"
$ws.Range("E2").Value = "94.68"
$ws.Range("F2").Value = "4.52"
$ws.Range("G2").Value = "0.6916"

$ws.Range("B3").Value = "009274"
$ws.Range("C3").Value = "融通健康产业灵活配置混合C"
$ws.Range("D3").Value = "3.16"
$ws.Range("E3").Value = "94.68"
$ws.Range("F3").Value = "4.52"
$ws.Range("G3").Value = "0.1428"

$ws.Range("B4").Value = "007808"
$ws.Range("C4").Value = "北信瑞丰量化优选灵活配置混合"
$ws.Range("D4").Value = "0.24"
$ws.Range("E4").Value = "79.84"
$ws.Range("F4").Value = "1.04"
$ws.Range("G4").Value = "0.0025"

$dataRange.ClearFormats()

# Column H ("仓位排名") is a plain number.
$ws.Range("H2").Value = 9
$ws.Range("H3").Value = 9
$ws.Range("H4").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" (grand-total) sheet: prepend a 2022-Q1 summary
#    row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2").EntireRow.Insert()

# Insert() can leave a stray blended style on the new row's cells; strip
# it back to the default (unstyled) look used by the other data rows.
$total.Range("A2:D2").ClearFormats()

# The index cell (A2) carries the same style as the rest of column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.84

# Renumber the 0-based index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
